$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 3172.4827  # was 3533.7666
$ws.Range("I92").Value = 1933.2727  # was 2458.3914
$ws.Range("K92").Value = 1933.2727  # was 2458.3914
$ws.Range("M92").Value = -685.2727  # was -1210.3914
$ws.Range("H135").Value = 1710.8889  # was 1653
$ws.Range("I135").Value = 1781  # was 1708.8889
$ws.Range("K135").Value = 16029  # was 15380.0001
$ws.Range("M135").Value = -13494  # was -12845.0001
$ws.Range("H137").Value = 3508.0908  # was 3609.45
$ws.Range("I137").Value = 1911.8  # was 1977.7142
$ws.Range("J137").Value = 6928.7144  # was 7416.8335
$ws.Range("K137").Value = 5735.4  # was 5933.142599999999
$ws.Range("L137").Value = 20786.1432  # was 22250.5005
$ws.Range("M137").Value = -3185.4  # was -3383.142599999999
$ws.Range("N137").Value = -25886.1432  # was -27350.5005
$ws.Range("H138").Value = 3860.926  # was 3864.611
$ws.Range("J138").Value = 4122.436  # was 4127.5386
$ws.Range("L138").Value = 12367.308  # was 12382.6158
$ws.Range("N138").Value = -22647.308  # was -22662.6158
$ws.Range("H141").Value = 2623.8333  # was 2636.375
$ws.Range("I141").Value = 2231.8096  # was 2308.45
$ws.Range("J141").Value = 5368  # was 4276
$ws.Range("K141").Value = 6695.4288  # was 6925.349999999999
$ws.Range("L141").Value = 16104  # was 12828
$ws.Range("M141").Value = -1515.4288  # was -1745.349999999999
$ws.Range("N141").Value = -26464  # was -23188

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6342.2964  # was 6142.5713
$ws.Range("I61").Value = 4706.174  # was 4541.3335
$ws.Range("K61").Value = 4706.174  # was 4541.3335
$ws.Range("M61").Value = -4494.174  # was -4329.3335
$ws.Range("H74").Value = 25649660  # was 33341496
$ws.Range("I74").Value = 30310596  # was 37044996
$ws.Range("J74").Value = 14507  # was 10000
$ws.Range("K74").Value = 30310596  # was 37044996
$ws.Range("L74").Value = 14507  # was 10000
$ws.Range("M74").Value = -30309722  # was -37044122
$ws.Range("N74").Value = -16255  # was -11748
$ws.Range("H76").Value = 30000  # was 0
$ws.Range("J76").Value = 30000  # was 0
$ws.Range("L76").Value = 30000  # was 0
$ws.Range("N76").Value = -30676  # was None
$ws.Range("H77").Value = 25649660  # was 33341496
$ws.Range("I77").Value = 30310596  # was 37044996
$ws.Range("J77").Value = 14507  # was 10000
$ws.Range("K77").Value = 151552980  # was 185224980
$ws.Range("L77").Value = 72535  # was 50000
$ws.Range("M77").Value = -151548612  # was -185220612
$ws.Range("N77").Value = -81271  # was -58736
$ws.Range("H79").Value = 30000  # was 0
$ws.Range("J79").Value = 30000  # was 0
$ws.Range("L79").Value = 30000  # was 0
$ws.Range("N79").Value = -32340  # was None
$ws.Range("H97").Value = 931.375  # was 943.8125
$ws.Range("I97").Value = 916.8  # was 930.06665
$ws.Range("K97").Value = 916.8  # was 930.06665
$ws.Range("M97").Value = -420.8  # was -434.06665
$ws.Range("H132").Value = 6519.1943  # was 5756.196
$ws.Range("I132").Value = 4795  # was 4401.143
$ws.Range("J132").Value = 17209.2  # was 10067.728
$ws.Range("K132").Value = 14385  # was 13203.429
$ws.Range("L132").Value = 51627.60000000001  # was 30203.184
$ws.Range("M132").Value = -11855  # was -10673.429
$ws.Range("N132").Value = -56687.60000000001  # was -35263.18399999999
$ws.Range("H136").Value = 6342.2964  # was 6142.5713
$ws.Range("I136").Value = 4706.174  # was 4541.3335
$ws.Range("K136").Value = 14118.522  # was 13624.0005
$ws.Range("M136").Value = -11568.522  # was -11074.0005
$ws.Range("H138").Value = 75429  # was 0
$ws.Range("J138").Value = 75429  # was 0
$ws.Range("L138").Value = 75429  # was 0
$ws.Range("N138").Value = -85709  # was None
$ws.Range("H139").Value = 79998  # was 0
$ws.Range("J139").Value = 79998  # was 0
$ws.Range("L139").Value = 79998  # was 0
$ws.Range("N139").Value = -90278  # was None

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1216.1428  # was 1336
$ws.Range("I94").Value = 878.25  # was 1005.3333
$ws.Range("K94").Value = 878.25  # was 1005.3333
$ws.Range("M94").Value = -427.25  # was -554.3333
$ws.Range("H132").Value = 69932.414  # was 69979.2
$ws.Range("J132").Value = 69932.414  # was 69979.2
$ws.Range("L132").Value = 69932.414  # was 69979.2
$ws.Range("N132").Value = -80052.414  # was -80099.2
$ws.Range("H134").Value = 2837.08  # was 2042.826
$ws.Range("I134").Value = 2163.0417  # was 1665.6888
$ws.Range("K134").Value = 6489.125100000001  # was 4997.0664
$ws.Range("M134").Value = -3954.125100000001  # was -2462.0664

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0  # was 7374.25
$ws.Range("I16").Value = 0  # was 5999.5
$ws.Range("J16").Value = 0  # was 8749
$ws.Range("K16").Value = 0  # was 5999.5
$ws.Range("L16").Value = 0  # was 8749
$ws.Range("M16").ClearContents()  # was -5712.5
$ws.Range("N16").ClearContents()  # was -9323
$ws.Range("H31").Value = 20267.5  # was 20890.662
$ws.Range("I31").Value = 3338  # was 3772.2727
$ws.Range("K31").Value = 3338  # was 3772.2727
$ws.Range("M31").Value = -3043  # was -3477.2727
$ws.Range("H34").Value = 20267.5  # was 20890.662
$ws.Range("I34").Value = 3338  # was 3772.2727
$ws.Range("K34").Value = 3338  # was 3772.2727
$ws.Range("M34").Value = -3136  # was -3570.2727
$ws.Range("H105").Value = 3991.8235  # was 4354.067
$ws.Range("I105").Value = 5110  # was 6250
$ws.Range("J105").Value = 3525.9167  # was 3664.6365
$ws.Range("K105").Value = 5110  # was 6250
$ws.Range("L105").Value = 3525.9167  # was 3664.6365
$ws.Range("M105").Value = -3363  # was -4503
$ws.Range("N105").Value = -7019.9167  # was -7158.636500000001
$ws.Range("H113").Value = 0  # was 7374.25
$ws.Range("I113").Value = 0  # was 5999.5
$ws.Range("J113").Value = 0  # was 8749
$ws.Range("K113").Value = 0  # was 5999.5
$ws.Range("L113").Value = 0  # was 8749
$ws.Range("M113").ClearContents()  # was -3829.5
$ws.Range("N113").ClearContents()  # was -13089
$ws.Range("H132").Value = 2656.8276  # was 2854.12
$ws.Range("I132").Value = 1478.8572  # was 1603.1111
$ws.Range("J132").Value = 5749  # was 6071
$ws.Range("K132").Value = 4436.571599999999  # was 4809.3333
$ws.Range("L132").Value = 17247  # was 18213
$ws.Range("M132").Value = -1906.571599999999  # was -2279.3333
$ws.Range("N132").Value = -22307  # was -23273
$ws.Range("H135").Value = 67626.10000000001  # was 68438.72
$ws.Range("J135").Value = 67626.10000000001  # was 68438.72
$ws.Range("L135").Value = 67626.10000000001  # was 68438.72
$ws.Range("N135").Value = -77766.10000000001  # was -78578.72

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 94302.64999999999  # was 89065.28
$ws.Range("J2").Value = 228895.14  # was 200287
$ws.Range("L2").Value = 1373370.84  # was 1201722
$ws.Range("N2").Value = -1373596.84  # was -1201948
$ws.Range("H4").Value = 2901589.5  # was 3445630.2
$ws.Range("I4").Value = 3286734.5  # was 4144133.2
$ws.Range("K4").Value = 9860203.5  # was 12432399.6
$ws.Range("M4").Value = -9860091.5  # was -12432287.6
$ws.Range("H6").Value = 62690.375  # was 50202.45
$ws.Range("I6").Value = 66869.60000000001  # was 55780.277
$ws.Range("K6").Value = 200608.8  # was 167340.831
$ws.Range("M6").Value = -200495.8  # was -167227.831
$ws.Range("H12").Value = 262.33334  # was 253.44
$ws.Range("I12").Value = 61.42857  # was 68.5
$ws.Range("J12").Value = 345.05884  # was 311.8421
$ws.Range("K12").Value = 184.28571  # was 205.5
$ws.Range("L12").Value = 1035.17652  # was 935.5263
$ws.Range("M12").Value = -11.28570999999999  # was -32.5
$ws.Range("N12").Value = -1381.17652  # was -1281.5263
$ws.Range("H17").Value = 221.44444  # was 209.3
$ws.Range("J17").Value = 325  # was 280
$ws.Range("L17").Value = 975  # was 840
$ws.Range("N17").Value = -1313  # was -1178
$ws.Range("H34").Value = 2356.8572  # was 2111.889
$ws.Range("J34").Value = 2399  # was 1826.75
$ws.Range("L34").Value = 7197  # was 5480.25
$ws.Range("N34").Value = -7365  # was -5648.25
$ws.Range("H39").Value = 4249.8335  # was 5139.9
$ws.Range("J39").Value = 4545.273  # was 5599.8887
$ws.Range("L39").Value = 13635.819  # was 16799.6661
$ws.Range("N39").Value = -14223.819  # was -17387.6661
$ws.Range("H55").Value = 1959.8  # was 1766.5
$ws.Range("J55").Value = 3447.5  # was 2565
$ws.Range("L55").Value = 10342.5  # was 7695
$ws.Range("N55").Value = -10696.5  # was -8049
$ws.Range("H117").Value = 4674.4165  # was 5036.4
$ws.Range("I117").Value = 2718.4285  # was 2966.6667
$ws.Range("J117").Value = 7412.8  # was 8141
$ws.Range("K117").Value = 8155.2855  # was 8900.000100000001
$ws.Range("L117").Value = 22238.4  # was 24423
$ws.Range("M117").Value = -4713.2855  # was -5458.000100000001
$ws.Range("N117").Value = -29122.4  # was -31307

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11205.104  # was 11363.857
$ws.Range("I122").Value = 11807.5625  # was 11871.125
$ws.Range("J122").Value = 10463.615  # was 10687.5
$ws.Range("K122").Value = 35422.6875  # was 35613.375
$ws.Range("L122").Value = 31390.845  # was 32062.5
$ws.Range("M122").Value = -32972.6875  # was -33163.375
$ws.Range("N122").Value = -36290.845  # was -36962.5
$ws.Range("H132").Value = 6902.6  # was 6585.3335
$ws.Range("I132").Value = 3874.75  # was 4099.6
$ws.Range("K132").Value = 11624.25  # was 12298.8
$ws.Range("M132").Value = -9094.25  # was -9768.800000000001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 27734.688  # was 28850.2
$ws.Range("I93").Value = 22169.285  # was 22740.428
$ws.Range("J93").Value = 32063.334  # was 34196.25
$ws.Range("K93").Value = 22169.285  # was 22740.428
$ws.Range("L93").Value = 32063.334  # was 34196.25
$ws.Range("M93").Value = -20921.285  # was -21492.428
$ws.Range("N93").Value = -34559.334  # was -36692.25
$ws.Range("H136").Value = 6742.5386  # was 7036.2705
$ws.Range("I136").Value = 4231.2085  # was 4496.909
$ws.Range("K136").Value = 12693.6255  # was 13490.727
$ws.Range("M136").Value = -10143.6255  # was -10940.727

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0  # was 478
$ws.Range("I8").Value = 0  # was 478
$ws.Range("K8").Value = 0  # was 478
$ws.Range("M8").ClearContents()  # was -338
$ws.Range("H96").Value = 1661  # was 1768.5714
$ws.Range("I96").Value = 1661  # was 1768.5714
$ws.Range("K96").Value = 1661  # was 1768.5714
$ws.Range("M96").Value = -288  # was -395.5714
$ws.Range("H136").Value = 3071.9119  # was 3430.4194
$ws.Range("I136").Value = 1113.3214  # was 1322.84
$ws.Range("K136").Value = 3339.9642  # was 3968.52
$ws.Range("M136").Value = -789.9642000000003  # was -1418.52
$ws.Range("H138").Value = 74293.336  # was 74195
$ws.Range("J138").Value = 74490  # was 74293.336
$ws.Range("L138").Value = 74490  # was 74293.336
$ws.Range("N138").Value = -84770  # was -84573.336
